# Edit script: insert "Overview" and "Game Demo" slides after the "Team MagiDev"
# slide, and bump the footer datetimeFigureOut field from 4/11/2017 to 4/12/2017
# across the slide master and all slide layouts.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1. Insert the new "Overview" slide at position 3 (Title and Content layout).
# ---------------------------------------------------------------------------
$overview = $p.Slides.Add(3, 2)

$titleShape = $overview.Shapes.Item(1)
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Text = "Overview"
$titleTr.ParagraphFormat.Alignment = 2

$bodyShape = $overview.Shapes.Item(2)
$bodyTr = $bodyShape.TextFrame.TextRange
$bodyTr.Text = "Game Demo`rNew Features`rNew Audio Assets`rSound Effect Showcase`rNew Art Assets`rArt Showcase`rBeta Targets`rFeature Content Target`rArt Content Target`rAudio Content Target`r"

# Paragraph-level formatting to match the target deck.
$bodyTr.Paragraphs(1,1).Runs(1,1).Font.Size = 24
$bodyTr.Paragraphs(2,1).Runs(1,1).Font.Size = 24
$bodyTr.Paragraphs(3,1).Runs(1,1).Font.Size = 24

$p4 = $bodyTr.Paragraphs(4,1)
$p4.IndentLevel = 3
$p4.Runs(1,1).Font.Size = 18

$bodyTr.Paragraphs(5,1).Runs(1,1).Font.Size = 24

$p6 = $bodyTr.Paragraphs(6,1)
$p6.IndentLevel = 3
$p6.Runs(1,1).Font.Size = 18

$bodyTr.Paragraphs(7,1).Runs(1,1).Font.Size = 24

$p8 = $bodyTr.Paragraphs(8,1)
$p8.IndentLevel = 3
$p8.Runs(1,1).Font.Size = 18

$p9 = $bodyTr.Paragraphs(9,1)
$p9.IndentLevel = 3
$p9.Runs(1,1).Font.Size = 18

$p10 = $bodyTr.Paragraphs(10,1)
$p10.IndentLevel = 3
$p10.Runs(1,1).Font.Size = 18

$p11 = $bodyTr.Paragraphs(11,1)
$p11.IndentLevel = 3
$p11.ParagraphFormat.Alignment = 2

# Corner logo picture, copied from the "Team MagiDev" slide.
$logoSource = $p.Slides.Item(2).Shapes.Item(3)
$logoSource.Copy()
$overview.Shapes.Paste() | Out-Null

# ---------------------------------------------------------------------------
# 2. Insert the new "Game Demo" slide at position 4 (Title and Content layout).
# ---------------------------------------------------------------------------
$gameDemo = $p.Slides.Add(4, 2)

$gdTitleShape = $gameDemo.Shapes.Item(1)
$gdTitleTr = $gdTitleShape.TextFrame.TextRange
$gdTitleTr.Text = "Game Demo"
$gdTitleTr.ParagraphFormat.Alignment = 2

# ---------------------------------------------------------------------------
# 3. Bump the datetimeFigureOut footer field on the master + all layouts.
# ---------------------------------------------------------------------------
function Update-DateShape($container) {
    for ($i = 1; $i -le $container.Shapes.Count; $i++) {
        $sh = $container.Shapes.Item($i)
        if ($sh.Type -eq 14 -and $sh.PlaceholderFormat.Type -eq 16) {
            if ($sh.TextFrame.TextRange.Text -eq "4/11/2017") {
                $sh.TextFrame.TextRange.Text = "4/12/2017"
            }
        }
    }
}

$master = $p.SlideMaster
Update-DateShape($master)

$layouts = $master.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DateShape($layouts.Item($j))
}
